$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new fiscal-year column (FY2018) was added to the report. This shifts the
# existing columns D:K one place to the right (to E:L) and the new column D
# is populated with the new year's figures.
$col = $ws.Range("D1").EntireColumn
$col.Insert()

# The freshly inserted column D picks up formatting from the neighboring
# column to its left by default; copy the (correct) number formats/styles
# back from column E (which now holds what used to be column D) so column D
# matches the rest of the data columns.
$ws.Range("E5:E102").Copy() | Out-Null
$ws.Range("D5:D102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the new column D with the FY2018 figures.
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 704600
$ws.Range("D9").Value2 = 489100
$ws.Range("D10").Value2 = 215500
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("D15").Value2 = 38800
$ws.Range("D17").Value2 = 603500
$ws.Range("D18").Value2 = 101000
$ws.Range("D20").Value2 = 1100
$ws.Range("D21").Value2 = 140900
$ws.Range("D22").Value2 = 8200
$ws.Range("D23").Value2 = 94000
$ws.Range("D24").Value2 = 20900
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 73100
$ws.Range("D27").Value2 = 72800
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = 100
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = -1100
$ws.Range("D33").Value2 = 72900
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 72900
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 23300
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 133400
$ws.Range("D44").Value2 = 0
$ws.Range("D45").Value2 = 21400
$ws.Range("D46").Value2 = 178100
$ws.Range("D47").Value2 = 0
$ws.Range("D48").Value2 = 52300
$ws.Range("D49").Value2 = 730200
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 7400
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 967900
$ws.Range("D57").Value2 = 16700
$ws.Range("D58").Value2 = 0
$ws.Range("D59").Value2 = 240300
$ws.Range("D60").Value2 = 257000
$ws.Range("D61").Value2 = 104200
$ws.Range("D62").Value2 = 17000
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 378200
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = -41500
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 589700
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 72900
$ws.Range("D83").Value2 = 38800
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 156600
$ws.Range("D91").Value2 = -16000
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -17000
$ws.Range("D96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = -141600
$ws.Range("D101").Value2 = -1200
$ws.Range("D102").Value2 = -3200

# "Capital Expenditures" row: besides the new FY2018 figure, the FY2017
# comparative figure (now in column E) was also restated.
$ws.Range("E91").Value2 = -11700
